# Insert a new data row at row 574 (pushing existing rows 574:622 down to
# 575:623) and populate it with a new price record, matching the commit
# "Fruta / hortaliza, semanal" weekly price-update pattern used throughout
# this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 574:622 down to 575:623, leaving a blank row 574 to fill in.
$ws.Rows.Item(574).Insert()

# Populate the newly inserted row 574 with the new record.
$ws.Cells.Item(574,1).Value = 3
$ws.Cells.Item(574,2).Value = "Femacal de La Calera"
$ws.Cells.Item(574,3).Value = "Coquimbo"
$ws.Cells.Item(574,4).Value = 45106
$ws.Cells.Item(574,5).Value = 5
$ws.Cells.Item(574,6).Value = 100112031
$ws.Cells.Item(574,7).Value = "Poroto verde"
$ws.Cells.Item(574,8).Value = "Magnum"
$ws.Cells.Item(574,9).Value = "Primera"
$ws.Cells.Item(574,10).Value = 85
$ws.Cells.Item(574,11).Value = 24000
$ws.Cells.Item(574,12).Value = 25000
$ws.Cells.Item(574,13).Value = 24471
$ws.Cells.Item(574,14).Value = "$/malla 25 kilos"
$ws.Cells.Item(574,15).Value = "Provincia de Limarí"
$ws.Cells.Item(574,16).Value = 979
$ws.Cells.Item(574,17).Value = 25
$ws.Cells.Item(574,18).Value = "Hortaliza"
